$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.48"
$ws.Range("E2").Value = "'0.26%"
$ws.Range("G2").Value = "'10"
$ws.Range("D3").Value = "'35.79"
$ws.Range("E3").Value = "'-0.32%"
$ws.Range("G3").Value = "'10"
$ws.Range("D4").Value = "'5.036"
$ws.Range("E4").Value = "'-1.16%"
$ws.Range("G4").Value = "'10"
$ws.Range("D5").Value = "'0.08042"
$ws.Range("E5").Value = "'-0.36%"
$ws.Range("G5").Value = "'10"
$ws.Range("D6").Value = "'1.885"
$ws.Range("E6").Value = "'-1.92%"
$ws.Range("G6").Value = "'10"
$ws.Range("D7").Value = "'4.147"
$ws.Range("E7").Value = "'-0.77%"
$ws.Range("G7").Value = "'10"
$ws.Range("D8").Value = "'7.781"
$ws.Range("E8").Value = "'0.49%"
$ws.Range("G8").Value = "'10"
$ws.Range("D9").Value = "'0.9204"
$ws.Range("E9").Value = "'-0.77%"
$ws.Range("G9").Value = "'10"
$ws.Range("D10").Value = "'0.1270"
$ws.Range("E10").Value = "'-5.21%"
$ws.Range("G10").Value = "'10"
$ws.Range("D11").Value = "'0.1914"
$ws.Range("E11").Value = "'0.50%"
$ws.Range("G11").Value = "'10"
$ws.Range("D12").Value = "'0.09077"
$ws.Range("E12").Value = "'-1.10%"
$ws.Range("G12").Value = "'10"
$ws.Range("D13").Value = "'0.03467"
$ws.Range("E13").Value = "'1.77%"
$ws.Range("G13").Value = "'10"
$ws.Range("E14").Value = "'0.30%"
$ws.Range("G14").Value = "'10"
$ws.Range("D15").Value = "'0.001414"
$ws.Range("E15").Value = "'0.21%"
$ws.Range("G15").Value = "'10"
$ws.Range("D16").Value = "'0.006209"
$ws.Range("E16").Value = "'5.00%"
$ws.Range("G16").Value = "'10"
$ws.Range("D17").Value = "'3.817"
$ws.Range("G17").Value = "'10"
$ws.Range("D18").Value = "'3.395"
$ws.Range("E18").Value = "'12.61%"
$ws.Range("G18").Value = "'10"
$ws.Range("D19").Value = "'0.3420"
$ws.Range("E19").Value = "'-0.94%"
$ws.Range("G19").Value = "'10"
$ws.Range("D20").Value = "'0.1321"
$ws.Range("E20").Value = "'-0.78%"
$ws.Range("G20").Value = "'10"
$ws.Range("D21").Value = "'5.193"
$ws.Range("E21").Value = "'5.86%"
$ws.Range("G21").Value = "'10"
$ws.Range("D22").Value = "'0.2307"
$ws.Range("G22").Value = "'10"
$ws.Range("D23").Value = "'0.04424"
$ws.Range("E23").Value = "'-0.35%"
$ws.Range("G23").Value = "'10"
$ws.Range("G24").Value = "'10"
$ws.Range("D25").Value = "'0.004610"
$ws.Range("E25").Value = "'-3.98%"
$ws.Range("G25").Value = "'10"
$ws.Range("G26").Value = "'10"
$ws.Range("E27").Value = "'-3.78%"
$ws.Range("G27").Value = "'10"
$ws.Range("E28").Value = "'41.82%"
$ws.Range("G28").Value = "'10"
$ws.Range("G29").Value = "'10"
$ws.Range("G30").Value = "'10"
$ws.Range("G31").Value = "'10"
$ws.Range("G32").Value = "'10"
$ws.Range("G33").Value = "'10"
$ws.Range("G34").Value = "'10"
$ws.Range("G35").Value = "'10"
$ws.Range("G36").Value = "'10"
$ws.Range("G37").Value = "'10"
$ws.Range("G38").Value = "'10"
$ws.Range("E39").Value = "'-2.45%"
$ws.Range("G39").Value = "'10"
$ws.Range("D40").Value = "'0.05350"
$ws.Range("E40").Value = "'8.99%"
$ws.Range("G40").Value = "'10"
$ws.Range("D41").Value = "'0.007616"
$ws.Range("E41").Value = "'-0.39%"
$ws.Range("G41").Value = "'10"
$ws.Range("D42").Value = "'0.01013"
$ws.Range("E42").Value = "'4.39%"
$ws.Range("G42").Value = "'10"
$ws.Range("D43").Value = "'0.1353"
$ws.Range("E43").Value = "'-1.66%"
$ws.Range("G43").Value = "'10"
$ws.Range("D44").Value = "'0.002153"
$ws.Range("E44").Value = "'2.30%"
$ws.Range("G44").Value = "'10"
$ws.Range("D45").Value = "'0.009604"
$ws.Range("E45").Value = "'-14.98%"
$ws.Range("G45").Value = "'10"
$ws.Range("D46").Value = "'0.00006126"
$ws.Range("E46").Value = "'-4.29%"
$ws.Range("G46").Value = "'10"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("G47").Value = "'10"
$ws.Range("E48").Value = "'2.60%"
$ws.Range("G48").Value = "'10"
$ws.Range("D49").Value = "'0.001661"
$ws.Range("E49").Value = "'39.18%"
$ws.Range("G49").Value = "'10"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("G50").Value = "'10"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.08%"
$ws.Range("G51").Value = "'10"

Write-Host "Updated 118 cells"
